$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 (shifts old rows 7-17 down to 8-18),
# inheriting styles/formatting from the surrounding rows.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new (unsuccessful) run's data.
$ws.Cells.Item(7,1).Value  = "ukb51139_subset.csv"
$ws.Cells.Item(7,2).Value  = "28012 x 1081"
$ws.Cells.Item(7,3).Value  = "all"
$ws.Cells.Item(7,4).Value  = "no events"
$ws.Cells.Item(7,5).Value  = "> 140/80"
$ws.Cells.Item(7,6).Value  = "zscore"
$ws.Cells.Item(7,7).Value  = "median"
$ws.Cells.Item(7,8).Value  = "age, sex"
$ws.Cells.Item(7,9).Value  = 25
$ws.Cells.Item(7,11).Value = 87
$ws.Cells.Item(7,12).Value = "89.7 & 83.3"
$ws.Cells.Item(7,13).Value = "81.6 & 70.3"
$ws.Cells.Item(7,14).Value = 18
$ws.Cells.Item(7,15).Value = 47.5081

# Row 7 has a slightly shorter custom row height than the rest.
$ws.Rows.Item(7).RowHeight = 18.75

# Mark cells A16:H16 (formerly the first blank "spacer" row) with a lone
# quote character, producing an empty, quote-prefixed string value.
$ws.Range("A16:H16").Value = "'"
